$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the EAN (column A) and image URL (column C) values between rows 32 and 33.
# (reading via .Value2 — .Value as a getter is unreliable in this host)
$a32 = $ws.Range("A32").Value2
$c32 = $ws.Range("C32").Value2
$a33 = $ws.Range("A33").Value2
$c33 = $ws.Range("C33").Value2

# Keep column A as text (it was stored as a string originally, e.g. "4003773034094")
# rather than letting Excel auto-coerce the EAN digits into a number.
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A33").NumberFormat = "@"

$ws.Range("A32").Value = $a33
$ws.Range("C32").Value = $c33
$ws.Range("A33").Value = $a32
$ws.Range("C33").Value = $c32
